$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02, 1.042221624808166, 1.044866057099057, 1.04930729370299, 1.057663511938765, 1.034452624761487, 1.047298612227613, 1.047635903941801, 1.052064716290133, 1.060397864627196, 1.019553204223243)
    3 = @(1.02, 1.043541132279711, 1.046100559512942, 1.050519502305762, 1.059052704395053, 1.034662365183026, 1.048262629974194, 1.048680628156228, 1.053088111783132, 1.06159947743304, 1.019881662141207)
    4 = @(1.02, 1.04439410479829, 1.046898827358242, 1.051303480226943, 1.059951466334157, 1.034796235155691, 1.048885116285145, 1.049355538934063, 1.053749350221759, 1.062376331971019, 1.020093545274233)
    5 = @(1.02, 1.044752498088162, 1.047234293549084, 1.051632971303466, 1.060329276395953, 1.034852072754924, 1.049146501547095, 1.04963901194066, 1.054027106147197, 1.062702764646174, 1.020182465867116)
    6 = @(1.02, 1.044812662501049, 1.047290612523476, 1.051688288924926, 1.060392710747471, 1.034861422271333, 1.049190371299885, 1.04968659317111, 1.054073729257289, 1.062757565018215, 1.020197386955313)
    7 = @(1.02, 1.044398894433511, 1.046903310363589, 1.051307883265461, 1.059956514761855, 1.034796982993812, 1.048888610134667, 1.049359327730365, 1.053753062505305, 1.062380694393156, 1.020094734043903)
    8 = @(1.02, 1.042667734284943, 1.045283376206215, 1.049717050357236, 1.058133026932216, 1.034523890144592, 1.04762467683475, 1.04798920194138, 1.052410779450964, 1.060804095848378, 1.019664343412386)
    9 = @(1.02, 1.039610616933114, 1.042424585732423, 1.046910583888675, 1.054918581912916, 1.034028498330026, 1.045387406215377, 1.045566347138422, 1.050037975950759, 1.058020655257326, 1.018900917345011)
    10 = @(1.02, 1.037567851137506, 1.040515662092358, 1.045037251671935, 1.052774559625035, 1.033688670824624, 1.043888959084115, 1.043945203943916, 1.048450868936881, 1.056161276149199, 1.018388538068452)
    11 = @(1.02, 1.036682141663203, 1.039688305351366, 1.044225476080716, 1.051845870293756, 1.033539243217778, 1.043238431785468, 1.043241789383875, 1.047762351899239, 1.055355207751762, 1.018165847187323)
    12 = @(1.02, 1.036352967640644, 1.039380867016289, 1.043923850932391, 1.05150086198052, 1.033483395710581, 1.042996539684025, 1.042980288621456, 1.0475064088047, 1.055055651965614, 1.018083004532336)
    13 = @(1.02, 1.036423584907043, 1.039446819131759, 1.043988554936428, 1.051574869851264, 1.033495390746217, 1.043048438059391, 1.043036391493903, 1.047561318448215, 1.055119914376061, 1.018100780252311)
    14 = @(1.02, 1.036654935779623, 1.039662894898337, 1.044200545623375, 1.051817352881488, 1.033534633856664, 1.043218442177587, 1.043220178184132, 1.047741199606665, 1.05533044938121, 1.018159001950796)
    15 = @(1.02, 1.036797454407994, 1.039796010129246, 1.044331147305482, 1.051966747683285, 1.033558767296259, 1.043323153128374, 1.043333385800497, 1.047852004087284, 1.055460147463962, 1.018194857619387)
    16 = @(1.02, 1.037626607500619, 1.04057055425667, 1.045091113356598, 1.052836186676073, 1.033698539728988, 1.043932096544475, 1.043991856419009, 1.0484965360827, 1.056214751993249, 1.018403299809003)
    17 = @(1.02, 1.038146394283163, 1.041056193633868, 1.045567653857534, 1.053381475575478, 1.033785604223071, 1.044313615841955, 1.044404507205307, 1.048900486334418, 1.056687839578956, 1.018533827858524)
    18 = @(1.02, 1.038449463763514, 1.04133938385899, 1.045845553238592, 1.053699503429454, 1.033836167538088, 1.044535986785207, 1.044645059878858, 1.049135979529863, 1.056963692668005, 1.018609882913832)
    19 = @(1.02, 1.03855278354902, 1.041435931804739, 1.045940299916801, 1.053807937771095, 1.033853371044252, 1.044611782018097, 1.044727058562453, 1.049216255671238, 1.057057736149139, 1.018635802227413)
    20 = @(1.02, 1.038090637843534, 1.041004096913252, 1.045516531665282, 1.053322974368541, 1.033776285782176, 1.044272699297926, 1.044360248142639, 1.04885715916361, 1.056637091191176, 1.018519831690833)
    21 = @(1.02, 1.036586813781084, 1.039599269348568, 1.044138122318442, 1.051745949138654, 1.033523087226565, 1.043168387332311, 1.043166063721932, 1.047688234561264, 1.055268456135745, 1.018141860580442)
    22 = @(1.02, 1.03564024388576, 1.038715295264308, 1.043270906911922, 1.050754109083341, 1.033361903757505, 1.042472572211683, 1.042413950550026, 1.046952142512219, 1.054407093647662, 1.017903489396212)
    23 = @(1.02, 1.036142140286336, 1.039183974646819, 1.043730688079898, 1.051279932416045, 1.033447538853376, 1.042841579389231, 1.042812782701733, 1.04734246833456, 1.054863799896536, 1.018029923618495)
    24 = @(1.02, 1.038115832128569, 1.041027637410409, 1.045539631764429, 1.053349408638077, 1.033780497064039, 1.044291188225831, 1.044380247338967, 1.04887673723055, 1.056660022483832, 1.018526156202926)
    25 = @(1.02, 1.040401762455423, 1.043164176831817, 1.04763652356533, 1.055749762394881, 1.034158252082702, 1.045967003789544, 1.046193741038861, 1.050652312818411, 1.058740887441604, 1.01909888178733)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
